# Apply lab 3 update: insert a new "stddev" column after TapPressure_H20 (now column C),
# shifting old FlowRate_GPM column from C to D, and append a new FlowRate_GPM_stddev
# column E with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column C (FlowRate_GPM), pushing it to column D.
$ws.Columns("C").Insert()

# New header + data for column C: TapPressure_H20_stddev
$ws.Range("C1").Value = "TapPressure_H20_stddev"
$ws.Range("C2").Value = 0.6026192827980198
$ws.Range("C3").Value = 0.66389005113799093
$ws.Range("C4").Value = 0.80729176881719866
$ws.Range("C5").Value = 0.66281973416608353
$ws.Range("C6").Value = 5.3743492629340732
$ws.Range("C7").Value = 4.2501729376579469
$ws.Range("C8").Value = 1.2570322191574885
$ws.Range("C9").Value = 0.82090194298710573
$ws.Range("C10").Value = 0.81452440110778523

# New header + data for column E: FlowRate_GPM_stddev
$ws.Range("E1").Value = "FlowRate_GPM_stddev"
$ws.Range("E2").Value = 0.11238327277669095
$ws.Range("E3").Value = 0.054037024344424693
$ws.Range("E4").Value = 0.060249481325569558
$ws.Range("E5").Value = 0.059749476985158451
$ws.Range("E6").Value = 0.043243496620878674
$ws.Range("E7").Value = 0.042071367935925842
$ws.Range("E8").Value = 0.14310835055998794
$ws.Range("E9").Value = 0.12557866060760406
$ws.Range("E10").Value = 0.078612976028133275

# Match formatting: columns C and E should match the centered style used across the sheet.
$ws.Range("C1:C10").HorizontalAlignment = -4108
$ws.Range("E1:E10").HorizontalAlignment = -4108

# Size the new columns to fit their (longer) header text, as Excel would after inserting
# the column and typing the new header (column D keeps the width it had as the original
# column C, so it is left untouched).
$ws.Columns("C").ColumnWidth = 20.91666667
$ws.Columns("E").ColumnWidth = 19.25

# Update the active selection to match the final state.
$ws.Range("D15").Select()
